# Generate Report for Handoff
# - Refresh the "Latest HO Xliff Generate Date" on the Overview sheet for the
#   rows that are "Ready for handoff" (also shared by de-de's "Latest
#   Handoff Datetime" text, which carried the identical timestamp string).
# - Refresh the "Latest Handoff Datetime" on the zh-cn sheet.
# - Mark the Priority column ("ht") for those same rows on both the zh-cn and
#   de-de language sheets.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 10, 12, 13, 14)

# Overview sheet: bump the handoff-generate timestamp.
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Cells.Item($r, 7).Value = "2016-08-26 08:24:59"
}

# zh-cn sheet: bump the handoff timestamp and set Priority to "ht".
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Cells.Item($r, 8).Value = "2016-08-26 08:24:54"
    $zhcn.Cells.Item($r, 5).Value = "ht"
}

# de-de sheet: bump the handoff timestamp (it mirrored the Overview sheet's
# value exactly) and set Priority to "ht".
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Cells.Item($r, 8).Value = "2016-08-26 08:24:59"
    $dede.Cells.Item($r, 5).Value = "ht"
}
